$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns I (I0) and J (IF)
# Copy the formatting (bold/border/center) from the existing H1 header cell
# first, then overwrite the values so the copied text doesn't stick.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I and J, one row per data row (rows 2-55)
$data = @(
    @(8,9),
    @(9,9),
    @(8,8),
    @(6,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(6,6),
    @(6,6),
    @(7,7),
    @(7,7),
    @(6,6),
    @(7,7),
    @(6,6),
    @(11,11),
    @(5,5),
    @(4,5),
    @(8,8),
    @(8,8),
    @(7,8),
    @(5,5),
    @(8,8),
    @(6,6),
    @(8,8),
    @(8,8),
    @(6,6),
    @(6,6),
    @(7,7),
    @(6,6),
    @(8,8),
    @(8,8),
    @(6,6),
    @(5,5),
    @(6,7),
    @(6,7),
    @(8,8),
    @(6,7),
    @(7,7),
    @(7,7),
    @(7,8),
    @(7,8),
    @(7,8),
    @(7,7),
    @(7,8),
    @(8,8),
    @(7,7),
    @(5,5),
    @(9,9),
    @(6,6),
    @(7,7),
    @(7,7),
    @(6,6),
    @(7,7)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
